# Applies the cryptos-list value updates captured in the commit diff.
# Values are written via a Formula(text-literal) -> Copy -> PasteSpecial(xlPasteValues)
# round-trip so that numeric-looking strings (e.g. "559.39") are stored back as
# plain text (matching the original inlineStr cells) instead of being auto-coerced
# into numbers by the normal Range.Value input-parsing path, and so that no new
# cell style/number-format entry gets minted as a side effect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($rangeAddr).Formula = '="' + $escaped + '"'
    $ws.Range($rangeAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

Set-TextValue 'D2' '59.257.00'
Set-TextValue 'E2' '  -2.94%  '

Set-TextValue 'D3' '2.584.21'
Set-TextValue 'E3' '  -2.53%  '

Set-TextValue 'E4' '  -0.05%  '

Set-TextValue 'D5' '559.39'
Set-TextValue 'E5' '  -2.21%  '

Set-TextValue 'D6' '143.21'
Set-TextValue 'E6' '  -2.72%  '

Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  +0.10%  '

Set-TextValue 'E8' '  -1.95%  '

Set-TextValue 'D9' '2.592.00'
Set-TextValue 'E9' '  -3.24%  '

Set-TextValue 'E10' '  -3.28%  '

Set-TextValue 'E11' '  -1.06%  '

Set-TextValue 'E12' '  +10.74%  '

Set-TextValue 'D13' '0.356'
Set-TextValue 'E13' '  +3.48%  '

Set-TextValue 'D14' '3.038.85'
Set-TextValue 'E14' '  -2.63%  '

Set-TextValue 'D15' '59.225.12'
Set-TextValue 'E15' '  -2.79%  '

Set-TextValue 'D16' '23.05'
Set-TextValue 'E16' '  +5.41%  '

Set-TextValue 'E17' '  -1.07%  '

Set-TextValue 'D18' '2.576.02'
Set-TextValue 'E18' '  -3.18%  '

Set-TextValue 'D19' '4.59'
Set-TextValue 'E19' '  +0.76%  '

Set-TextValue 'D20' '336.85'
Set-TextValue 'E20' '  -2.36%  '

Set-TextValue 'D21' '10.37'
Set-TextValue 'E21' '  -1.17%  '

Set-TextValue 'E22' '  -0.06%  '

Set-TextValue 'E23' '  +0.19%  '

Set-TextValue 'D24' '64.03'
Set-TextValue 'E24' '  -4.39%  '

Set-TextValue 'D25' '0.467'
Set-TextValue 'E25' '  +5.18%  '

Set-TextValue 'D26' '0.999'

Set-TextValue 'E27' '  -2.78%  '

Set-TextValue 'D28' '7.36'
Set-TextValue 'E28' '  -0.75%  '

Set-TextValue 'D29' '0.0₃0776'
Set-TextValue 'E29' '  -1.72%  '

Set-TextValue 'E30' '  +0.03%  '

Set-TextValue 'E31' '  +0.00%  '

Set-TextValue 'E32' '  -3.37%  '

Set-TextValue 'D33' '159.11'
Set-TextValue 'E33' '  +2.73%  '

Set-TextValue 'D34' '19.03'
Set-TextValue 'E34' '  -1.36%  '

Set-TextValue 'D35' '4.04'
Set-TextValue 'E35' '  -1.80%  '

Set-TextValue 'E36' '  -1.80%  '

Set-TextValue 'D37' '0.881'
Set-TextValue 'E37' '  -3.91%  '

Set-TextValue 'D38' '0.867'
Set-TextValue 'E38' '  -5.06%  '

Set-TextValue 'E39' '  -0.57%  '

Set-TextValue 'E40' '  -2.65%  '

Set-TextValue 'D41' '3.68'
Set-TextValue 'E41' '  +0.22%  '

Set-TextValue 'D42' '292.42'
Set-TextValue 'E42' '  -5.07%  '

Set-TextValue 'D43' '132.49'
Set-TextValue 'E43' '  +4.00%  '

Set-TextValue 'E44' '  +0.24%  '

Set-TextValue 'D45' '0.0975'
Set-TextValue 'E45' '  -0.82%  '

Set-TextValue 'D46' '0.597'
Set-TextValue 'E46' '  -2.06%  '

Set-TextValue 'B47' 'Hedera'
Set-TextValue 'C47' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D47' '0.0535'
Set-TextValue 'E47' '  -2.85%  '

Set-TextValue 'B48' 'WhiteBITCoin'
Set-TextValue 'C48' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D48' '10.63'
Set-TextValue 'E48' '  -0.27%  '

Set-TextValue 'D49' '0.0234'
Set-TextValue 'E49' '  -0.51%  '

Set-TextValue 'B50' 'InjectiveProtocol'
Set-TextValue 'C50' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D50' '18.65'
Set-TextValue 'E50' '  -1.49%  '

Set-TextValue 'B51' 'Maker'
Set-TextValue 'C51' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D51' '1.955.00'
Set-TextValue 'E51' '  -0.67%  '
